$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (right after the header row), shifting the
# existing data rows down by one.
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new car entry.
$ws.Range("A2").Value = "Honda Amaze 2013-2016 Honda Amaze VX AT i-Vtech"
$ws.Range("B2").Value = "Rs. 3.90 Lakh"

# The insert pushed every old row down by one. The old row 12 (a duplicate
# "Honda Civic" entry identical to the row above it) is now at row 13;
# remove it so the table stays at 13 rows total and the final row is the
# original last row ("Honda Brio ... E MT").
$ws.Rows("13:13").Delete()
